$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove C2 entirely (was 0.7771450785698075)
$ws.Range("C2").ClearContents()

# Update remaining values with their new, slightly adjusted values
$ws.Range("E2").Value = 0.4748521911469572
$ws.Range("C4").Value = 0.5978820435290855
$ws.Range("C5").Value = -0.50613598754502
$ws.Range("C6").Value = -0.2706540469742613
$ws.Range("E7").Value = 0.01247916696662799
$ws.Range("E8").Value = 0.174086048246691
$ws.Range("E9").Value = -0.03768624985648339
$ws.Range("C10").Value = -0.03047919532178645
$ws.Range("E10").Value = -0.1249617237519041
$ws.Range("C11").Value = 0.1932702877606163
$ws.Range("E11").Value = -0.2500935825088479
$ws.Range("C14").Value = -0.2932081122163255
$ws.Range("E14").Value = -0.112644651861793
$ws.Range("C15").Value = -0.02187747290984809
$ws.Range("E15").Value = -0.1415666278731686
$ws.Range("E16").Value = 0.9379151023484189
$ws.Range("E19").Value = -0.338776212162295
